$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header cells: "_old" columns become "_FV2310", "_new" columns
#    become "_FV2404" (the "diff" column name is unchanged).
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$headerRange = $ws.Range("A1:U1")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Remember / reassert the original header formatting (bold font, grey fill,
# thin border, centered + wrapped text) so that turning the range into a
# table below does not bake the formatting into a table "headerRowDxfId".
$headerRange.ClearFormats()

# ---------------------------------------------------------------------------
# 2) Turn A1:U85 into an Excel Table ("Table1") with the header row as
#    column headers.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U85")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Re-apply the header formatting directly on the cells (matches the
# pre-existing style used by row 1 before the edit).
$headerRange.Font.Bold = $true
$headerRange.Interior.Pattern = 1
$headerRange.Interior.Color = 14277081
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split at row 2, top-left cell A2, frozen pane).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
